# Refresh the crypto price/volume table (and a few rank-order swaps)
# to match the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.181.65"
$ws.Range("E2").Value = "  +1.35%  "

$ws.Range("D3").Value = "2.245.69"
$ws.Range("E3").Value = "  +0.96%  "

$ws.Range("D4").Value = "'1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "'306.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.14%  "

$ws.Range("D6").Value = "'95.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.56%  "

$ws.Range("D7").Value = "'0.574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.85%  "

$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").Value = "'0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.29%  "

$ws.Range("D10").Value = "'35.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.24%  "

$ws.Range("D11").Value = "'0.0815"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("D12").Value = "'7.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.59%  "

$ws.Range("D13").Value = "'0.104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.17%  "

$ws.Range("D14").Value = "2.588.41"
$ws.Range("E14").Value = "  +1.01%  "

$ws.Range("D15").Value = "2.319.88"
$ws.Range("E15").Value = "  +4.05%  "

$ws.Range("D16").Value = "'0.835"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.38%  "

$ws.Range("D17").Value = "'13.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.15%  "

$ws.Range("D18").Value = "44.033.17"
$ws.Range("E18").Value = "  +1.25%  "

$ws.Range("D19").Value = "0.0₃0972"
$ws.Range("E19").Value = "  +1.04%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.61%  "

$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "'12.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.26%  "

$ws.Range("D22").Value = "'65.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.72%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'236.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.00%  "

$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'2.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.65%  "

$ws.Range("E25").Value = "  -0.99%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").Value = "'10.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.15%  "

$ws.Range("E28").Value = "  +1.13%  "

$ws.Range("D29").Value = "'37.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.91%  "

$ws.Range("D30").Value = "'6.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.08%  "

$ws.Range("D31").Value = "'20.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.32%  "

$ws.Range("D32").Value = "'153.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.36%  "

$ws.Range("D33").Value = "'0.0803"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.07%  "

$ws.Range("D34").Value = "'3.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.18%  "

$ws.Range("D35").Value = "'2.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.82%  "

$ws.Range("E36").Value = "  +2.83%  "

$ws.Range("D37").Value = "'0.108"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").Value = "'1.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.79%  "

$ws.Range("E39").Value = "  -2.68%  "

$ws.Range("D40").Value = "'3.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.80%  "

$ws.Range("D41").Value = "'14.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.15%  "

$ws.Range("D42").Value = "'0.0298"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.85%  "

$ws.Range("E43").Value = "  +0.23%  "

$ws.Range("D44").Value = "1.737.70"
$ws.Range("E44").Value = "  +1.84%  "

$ws.Range("D45").Value = "'83.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.97%  "

$ws.Range("D46").Value = "'0.191"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.39%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'14.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.50%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'100.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.53%  "

$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'4.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.71%  "

$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'8.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.71%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'54.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.78%  "
